# Semana 6, Clase 12
# Fills in the Solver results for the "Ejemplo Lineal" and "Ejemplo No Lineal"
# worksheets (values + formulas that were previously blank placeholders),
# flips the "Ejemplo No Lineal" objective from MIN(-XY) to MAX(XY) with an
# equality constraint, updates the related Solver-engine defined names, and
# makes "Ejemplo Lineal" the active sheet (was "Ejemplo Sistema Bombeo").

$wb = $excel.ActiveWorkbook

$wsLineal    = $wb.Worksheets.Item("Ejemplo Lineal")
$wsNoLineal  = $wb.Worksheets.Item("Ejemplo No Lineal")

# ---------------------------------------------------------------------
# Sheet "Ejemplo Lineal": fill in the optimal X/Y plus every dependent
# formula cell that Solver would have populated.
# ---------------------------------------------------------------------
$wsLineal.Range("E2").Value = 7.7272727272727284
$wsLineal.Range("E3").Value = 4.545454545454545

$wsLineal.Range("E6").Formula = "=E2+2*E3"

$wsLineal.Range("C9").Formula  = "=2*E2+E3"
$wsLineal.Range("E9").Value    = 20

$wsLineal.Range("C10").Formula = "=-4*E2+5*E3"
$wsLineal.Range("E10").Value   = 10

$wsLineal.Range("C11").Formula = "=-E2+2*E3"
$wsLineal.Range("E11").Value   = -2

$wsLineal.Range("C12").Formula = "=-E2+5*E3"
$wsLineal.Range("E12").Value   = 15

$wsLineal.Range("C13").Formula = "=E2"
$wsLineal.Range("E13").Value   = 0

$wsLineal.Range("C14").Formula = "=E3"
$wsLineal.Range("E14").Value   = 0

# ---------------------------------------------------------------------
# Sheet "Ejemplo No Lineal": fill in the optimal X/Y, switch the
# objective from MIN( -XY) to MAX(XY), change Restr. 1 to an equality
# constraint, and fill in the dependent formula cells.
# ---------------------------------------------------------------------
$wsNoLineal.Range("E2").Value = 12.000000249999999
$wsNoLineal.Range("E3").Value = 12.000000249999999

$wsNoLineal.Range("C6").Value = "MAX"
$wsNoLineal.Range("D6").Value = "XY"
$wsNoLineal.Range("E6").Formula = "=E2*E3"

$wsNoLineal.Range("C9").Formula = "=2*E2+2*E3"
# D9 needs the literal text "=" (not a formula) while keeping its style;
# assigning "=" directly is parsed as a formula, so build it as a text
# formula first and then convert that cell to a plain value in place.
$wsNoLineal.Range("D9").Formula = "=""="""
$wsNoLineal.Range("D9").Copy() | Out-Null
$wsNoLineal.Range("D9").PasteSpecial(-4163) | Out-Null
$wsNoLineal.Range("E9").Value = 48

$wsNoLineal.Range("C10").Formula = "=E2"
$wsNoLineal.Range("E10").Value = 0

$wsNoLineal.Range("C11").Formula = "=E3"
$wsNoLineal.Range("E11").Value = 0

# ---------------------------------------------------------------------
# Solver defined names: point solver_adj / solver_opt at the cells that
# now hold the decision variables / objective, refresh solver_num and
# solver_typ (Solver now has the right-hand/engine bookkeeping for a
# resolved model), for both the linear and non-linear sheets.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    switch ($n.Name) {
        "Ejemplo Lineal!solver_adj"   { $n.RefersTo = "='Ejemplo Lineal'!`$E`$2:`$E`$3" }
        "Ejemplo No Lineal!solver_adj" { $n.RefersTo = "='Ejemplo No Lineal'!`$E`$2:`$E`$3" }
        "Ejemplo Lineal!solver_opt"   { $n.RefersTo = "='Ejemplo Lineal'!`$E`$6" }
        "Ejemplo No Lineal!solver_opt" { $n.RefersTo = "='Ejemplo No Lineal'!`$E`$6" }
        "Ejemplo Lineal!solver_num"   { $n.RefersTo = "=6" }
        "Ejemplo No Lineal!solver_num" { $n.RefersTo = "=3" }
        "Ejemplo Lineal!solver_typ"   { $n.RefersTo = "=1" }
        "Ejemplo No Lineal!solver_typ" { $n.RefersTo = "=1" }
    }
}

# ---------------------------------------------------------------------
# Active sheet moves from "Ejemplo Sistema Bombeo" to "Ejemplo Lineal".
# ---------------------------------------------------------------------
$wsLineal.Activate()
